$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = 1.36

# Row 7 updates
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.85
